$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value  = 0
$ws.Range("C3").Value  = 0.3677972307464058
$ws.Range("C4").Value  = 0.02992461252661037
$ws.Range("C5").Value  = 0.1234796466906095
$ws.Range("C6").Value  = 0.2978233793571007
$ws.Range("C7").Value  = 0.01090413015614074
$ws.Range("C8").Value  = 0
$ws.Range("C9").Value  = 0.04514393665358528
$ws.Range("C10").Value = 0.04884331010259938
$ws.Range("C11").Value = 0.04509121361603211
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0.03099254015091614
